# Fixed data conversion error for blank values passed in excel
# Row 2 originally held a fully-populated sample row (Book / Reading / 1 /
# 05/05/2024 / 05/10/2024). The fix replaces it with a regression-test row
# that exercises blank CreatedDate + Description handling: only Name (A2),
# Quantity (C2) and UpdatedDate (E2) stay populated; Description (B2) and
# CreatedDate (D2) are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values -----------------------------------------------------
$ws.Range("A2").Value2 = "Test25- Blank Created Date and description"
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value2 = 1
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value2 = "05/05/2024"

# --- Column widths (best-fit to new content) ---------------------------
$ws.Columns.Item(1).ColumnWidth = 22.6666666666667
$ws.Columns.Item(2).ColumnWidth = 26.6666666666667
$ws.Columns.Item(5).ColumnWidth = 11

# --- Selection moves to C4 after the edit -------------------------------
$ws.Range("C4").Select()
